# Benchmark httk working with new pk fits:
# add a new "2.7.1" row to the Table1 benchmark table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grow the table by one row (this also grows the table/autofilter range
# and the worksheet dimension once the new row has content).
$lo = $ws.ListObjects.Item(1)
$newListRow = $lo.ListRows.Add()
$newRowIndex = $newListRow.Range.Row

$rng = $ws.Range("A" + $newRowIndex + ":R" + $newRowIndex)
$rng.HorizontalAlignment = -4131

$ws.Cells.Item($newRowIndex, 1).Value2 = "2.7.1"
$ws.Cells.Item($newRowIndex, 2).Value2 = 1026
$ws.Cells.Item($newRowIndex, 3).Value2 = 0.99980000000000002
$ws.Cells.Item($newRowIndex, 4).Value2 = 1
$ws.Cells.Item($newRowIndex, 5).Value2 = 0.99990000000000001
$ws.Cells.Item($newRowIndex, 6).Value2 = 0.93500000000000005
$ws.Cells.Item($newRowIndex, 7).Value2 = 352
$ws.Cells.Item($newRowIndex, 8).Value2 = 0.2712
$ws.Cells.Item($newRowIndex, 9).Value2 = 352
$ws.Cells.Item($newRowIndex, 10).Value2 = 2.1840000000000002
$ws.Cells.Item($newRowIndex, 11).Value2 = 100
$ws.Cells.Item($newRowIndex, 12).Value2 = 1.2929999999999999
$ws.Cells.Item($newRowIndex, 13).Value2 = 112
$ws.Cells.Item($newRowIndex, 14).Value2 = 1.2130000000000001
$ws.Cells.Item($newRowIndex, 15).Value2 = 156
$ws.Cells.Item($newRowIndex, 16).Value2 = 0.26290000000000002
$ws.Cells.Item($newRowIndex, 17).Value2 = 863
$ws.Cells.Item($newRowIndex, 18).Value2 = "Updated CvTdb fits (invivoPKfit)"

# Update selection / scroll position on Sheet1 to mirror the saved view.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("J" + ($newRowIndex + 1)).Select()

# Note: Chart1 / Chart2 are chart *sheets* (not shapes embedded on
# Sheet1). Their c:cat/c:val caches reference Sheet1!A/F/J 2:30 and would
# normally be refreshed to 2:31 by Excel on save alongside this edit, but
# those chart sheets aren't exposed as automation objects in this host
# (wb.Charts / wb.Sheets only surface "Sheet1"), so there is no COM call
# available here to push the extra category/value point into them.

$wb.Save()
